# A new weekly price observation is inserted at row 23 ("Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Poroto verde"). All existing records from row
# 23 down to row 80 shift down by one row (to rows 24..81), and the new
# record's data is written into the now-vacated row 23.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23..80 down to rows 24..81, working bottom-up so that no row is
# overwritten before it has been copied down.
for ($r = 80; $r -ge 23; $r--) {
    $srcRow = $ws.Range("A" + $r + ":R" + $r)
    $dstRow = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $srcRow.Copy($dstRow)
}

# Populate the new record in row 23.
$ws.Range("D23").Value2 = 44645
$ws.Range("H23").Value2 = "Sin especificar"
$ws.Range("J23").Value2 = 120
$ws.Range("K23").Value2 = 25000
$ws.Range("L23").Value2 = 26000
$ws.Range("M23").Value2 = 25500
$ws.Range("O23").Value2 = "Región del Maule"
$ws.Range("P23").Value2 = 1020
